$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# These Price-column values look like plain numbers (e.g. "1.00", "7.17").
# The source stores them as literal text, so force a Text number format
# on just those specific cells before assigning, preventing Excel from
# auto-converting them to numeric values (which would drop formatting like
# trailing zeros, e.g. "1.00" -> 1).
$textCells = "D4","D5","D6","D9","D10","D11","D18","D19","D21","D22","D23","D24","D25","D27","D28","D32","D34","D35","D36","D41","D43","D45","D48","D49","D50"
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = '71.189.20'
$ws.Range("E2").Value = '  +0.49%  '
$ws.Range("D3").Value = '3.815.33'
$ws.Range("E3").Value = '  -0.76%  '
$ws.Range("D4").Value = '1.00'
$ws.Range("E4").Value = '  +0.01%  '
$ws.Range("D5").Value = '705.60'
$ws.Range("E5").Value = '  +1.02%  '
$ws.Range("D6").Value = '172.31'
$ws.Range("E6").Value = '  +0.28%  '
$ws.Range("D7").Value = '3.814.81'
$ws.Range("E7").Value = '  -0.73%  '
$ws.Range("E8").Value = '  +0.03%  '
$ws.Range("D9").Value = '0.523'
$ws.Range("E9").Value = '  +0.03%  '
$ws.Range("D10").Value = '0.161'
$ws.Range("E10").Value = '  -0.15%  '
$ws.Range("D11").Value = '7.63'
$ws.Range("E11").Value = '  +5.17%  '
$ws.Range("E12").Value = '  +0.82%  '
$ws.Range("E13").Value = '  -1.09%  '
$ws.Range("E14").Value = '  -0.43%  '
$ws.Range("D15").Value = '4.460.17'
$ws.Range("E15").Value = '  -0.73%  '
$ws.Range("D16").Value = '3.825.64'
$ws.Range("E16").Value = '  -0.46%  '
$ws.Range("D17").Value = '71.208.91'
$ws.Range("E17").Value = '  +0.44%  '
$ws.Range("B18").Value = 'Polkadot'
$ws.Range("C18").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range("D18").Value = '7.17'
$ws.Range("E18").Value = '  +0.18%  '
$ws.Range("B19").Value = 'Chainlink'
$ws.Range("C19").Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range("D19").Value = '17.52'
$ws.Range("E19").Value = '  +0.79%  '
$ws.Range("E20").Value = '  -0.47%  '
$ws.Range("D21").Value = '511.79'
$ws.Range("E21").Value = '  +3.56%  '
$ws.Range("D22").Value = '10.69'
$ws.Range("E22").Value = '  -0.39%  '
$ws.Range("D23").Value = '0.723'
$ws.Range("E23").Value = '  +0.85%  '
$ws.Range("D24").Value = '84.16'
$ws.Range("D25").Value = '0.0000144'
$ws.Range("E25").Value = '  -1.77%  '
$ws.Range("D26").Value = '3.966.88'
$ws.Range("E26").Value = '  -0.79%  '
$ws.Range("D27").Value = '12.06'
$ws.Range("E27").Value = '  -0.74%  '
$ws.Range("D28").Value = '10.41'
$ws.Range("E28").Value = '  -1.23%  '
$ws.Range("E29").Value = '  +0.16%  '
$ws.Range("E30").Value = '  -3.46%  '
$ws.Range("E31").Value = '  -2.83%  '
$ws.Range("D32").Value = '7.43'
$ws.Range("E32").Value = '  -0.93%  '
$ws.Range("E33").Value = '  -0.44%  '
$ws.Range("D34").Value = '29.12'
$ws.Range("E34").Value = '  -0.96%  '
$ws.Range("D35").Value = '0.173'
$ws.Range("E35").Value = '  -5.17%  '
$ws.Range("D36").Value = '9.19'
$ws.Range("E36").Value = '  +0.58%  '
$ws.Range("D37").Value = '3.777.72'
$ws.Range("E37").Value = '  -0.62%  '
$ws.Range("E38").Value = '  -0.08%  '
$ws.Range("E39").Value = '  -1.69%  '
$ws.Range("E40").Value = '  -0.11%  '
$ws.Range("D41").Value = '6.02'
$ws.Range("E41").Value = '  +0.45%  '
$ws.Range("E42").Value = '  -1.37%  '
$ws.Range("D43").Value = '3.30'
$ws.Range("E43").Value = '  -0.44%  '
$ws.Range("E44").Value = '  -0.02%  '
$ws.Range("D45").Value = '170.75'
$ws.Range("E45").Value = '  +4.45%  '
$ws.Range("E46").Value = '  +0.06%  '
$ws.Range("E47").Value = '  -0.25%  '
$ws.Range("D48").Value = '49.68'
$ws.Range("E48").Value = '  +1.98%  '
$ws.Range("D49").Value = '427.46'
$ws.Range("E49").Value = '  +5.76%  '
$ws.Range("D50").Value = '8.69'
$ws.Range("E50").Value = '  +0.79%  '
$ws.Range("E51").Value = '  -1.40%  '
